# Applies the "1401ME19" marksheet correction:
#  - recompute Right/Wrong/NotAttempt/Max + Marking + Total numbers
#  - drop the unused third Student/Correct-Ans column-pair (G:H)
#  - backfill the "Student Ans" column (A) -- and, for the first few rows,
#    column D -- with the correct-answer text, colour-coded green/red
#    depending on whether it agrees with the already-present answer

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Summary block (rows 10-12)
# ---------------------------------------------------------------------

# Row headers (No./Marking/Total in col A) pick up the "mtitleStyle" look
# already used by the row-9 header cells.
$ws.Range("A9").Copy()
$ws.Range("A10:A12").PasteSpecial(-4122)

$ws.Range("B10").Value = 19
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 8
$ws.Range("E10").Value = 28

$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

$ws.Range("B12").Value = 76
$ws.Range("C12").Value = -1
$ws.Range("E12").Value = "75/112"

# ---------------------------------------------------------------------
# Backfill column D for the first three answer rows (16-18)
# ---------------------------------------------------------------------

$ws.Range("B10").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("D16").Value = "Option A"

$ws.Range("B10").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("D17").Value = "Option C"

$ws.Range("B10").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("D18").Value = "Option D"

# ---------------------------------------------------------------------
# Backfill column A ("Student Ans") with the correct-answer text for
# every answered row, using the green "correctStyle" when it matches the
# student's answer in column B and the red "incorrectStyle" otherwise.
# ---------------------------------------------------------------------

function Set-AnswerCell($row, $style, $text) {
    if ($style -eq "correct") {
        $ws.Range("B10").Copy()
    } else {
        $ws.Range("C10").Copy()
    }
    $cell = $ws.Range("A$row")
    $cell.PasteSpecial(-4122)
    $cell.Value = $text
}

Set-AnswerCell 16 "correct" "Option A"
Set-AnswerCell 18 "correct" "Option B"
Set-AnswerCell 19 "correct" "Option C"
Set-AnswerCell 21 "correct" "Option C"
Set-AnswerCell 22 "correct" "Option D"
Set-AnswerCell 23 "correct" "Option D"
Set-AnswerCell 25 "correct" "Option A"
Set-AnswerCell 27 "correct" "Option A"
Set-AnswerCell 28 "correct" "Option D"
Set-AnswerCell 29 "incorrect" "Option C"
Set-AnswerCell 30 "correct" "Option B"
Set-AnswerCell 32 "correct" "Option C"
Set-AnswerCell 33 "correct" "Option D"
Set-AnswerCell 36 "correct" "Option A"
Set-AnswerCell 38 "correct" "Option A"
Set-AnswerCell 39 "correct" "Option D"
Set-AnswerCell 40 "correct" "Option D"

# ---------------------------------------------------------------------
# The second Student-Ans/Correct-Ans column pair (D:E) is only kept for
# the first three answer rows; clear it out everywhere else.
# ---------------------------------------------------------------------

$ws.Range("D19:E40").Clear()

# ---------------------------------------------------------------------
# Drop the unused third Student-Ans/Correct-Ans column pair entirely.
# ---------------------------------------------------------------------

$ws.Range("G1:H1").EntireColumn.Delete()
